$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The table now reports a "mean" and "std" column for every horizon,
#     the algorithm list changed (CART -> DTREE, NB removed), and the
#     result set only has 7 algorithms (one fewer row). Rewrite the
#     whole block with the new data. ---

# Drop the old last data row (previously algorithm index 7 / "NB");
# the refreshed results only include 7 algorithms.
$ws.Rows.Item(9).Delete()

# Header row: B1 keeps "Algorithm"; C1:L1 become the new mean/std pairs.
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Base mean"
$ws.Range("D1").Value = "One Year Base std"
$ws.Range("E1").Value = "Two Year Base mean"
$ws.Range("F1").Value = "Two Year Base std"
$ws.Range("G1").Value = "Three Year Base mean"
$ws.Range("H1").Value = "Three Year Base std"
$ws.Range("I1").Value = "Five Year Base mean"
$ws.Range("J1").Value = "Five Year Base std"
$ws.Range("K1").Value = "Ten Year Base mean"
$ws.Range("L1").Value = "Ten Year Base std"

# Row 2 - LR
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8647381650607457
$ws.Range("D2").Value = 0.02143645563772809
$ws.Range("E2").Value = 0.8515903658760801
$ws.Range("F2").Value = 0.03049269606396709
$ws.Range("G2").Value = 0.8325591848966137
$ws.Range("H2").Value = 0.0368960374448482
$ws.Range("I2").Value = 0.8250847868217054
$ws.Range("J2").Value = 0.05002388606835779
$ws.Range("K2").Value = 0.8389507318078747
$ws.Range("L2").Value = 0.04306738218990206

# Row 3 - LDA
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.8375576036866359
$ws.Range("D3").Value = 0.01887976606521018
$ws.Range("E3").Value = 0.8225087332230189
$ws.Range("F3").Value = 0.03592383650015361
$ws.Range("G3").Value = 0.8113525122365397
$ws.Range("H3").Value = 0.02862950412618045
$ws.Range("I3").Value = 0.802531492248062
$ws.Range("J3").Value = 0.05449779360672961
$ws.Range("K3").Value = 0.8065038136466708
$ws.Range("L3").Value = 0.04520591062944009

# Row 4 - KNN
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.8873900293255131
$ws.Range("D4").Value = 0.03012607663383144
$ws.Range("E4").Value = 0.8983590733590733
$ws.Range("F4").Value = 0.02627559412815156
$ws.Range("G4").Value = 0.8954250324642892
$ws.Range("H4").Value = 0.02785756923757248
$ws.Range("I4").Value = 0.8950823643410853
$ws.Range("J4").Value = 0.03036409627737291
$ws.Range("K4").Value = 0.8844258915687486
$ws.Range("L4").Value = 0.04179079833207544

# Row 5 - DTREE (was CART)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.7916547968160872
$ws.Range("D5").Value = 0.03783201731496429
$ws.Range("E5").Value = 0.7872448979591838
$ws.Range("F5").Value = 0.04141869901784344
$ws.Range("G5").Value = 0.7745729697332934
$ws.Range("H5").Value = 0.02698937209257003
$ws.Range("I5").Value = 0.7651647286821704
$ws.Range("J5").Value = 0.035318261666922
$ws.Range("K5").Value = 0.7842094413522985
$ws.Range("L5").Value = 0.04484765879048858

# Row 6 - RTREE
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8537243401759531
$ws.Range("D6").Value = 0.02083139725462477
$ws.Range("E6").Value = 0.8455460562603421
$ws.Range("F6").Value = 0.03463259763859165
$ws.Range("G6").Value = 0.8403256417940266
$ws.Range("H6").Value = 0.02482491583317106
$ws.Range("I6").Value = 0.8242974806201551
$ws.Range("J6").Value = 0.04788590580191691
$ws.Range("K6").Value = 0.8115852401566688
$ws.Range("L6").Value = 0.04440990971315106

# Row 7 - XTREE
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8336866359447006
$ws.Range("D7").Value = 0.02215932895701345
$ws.Range("E7").Value = 0.8346892811178526
$ws.Range("F7").Value = 0.04525953605319261
$ws.Range("G7").Value = 0.8304365198281889
$ws.Range("H7").Value = 0.03262255208120857
$ws.Range("I7").Value = 0.8374939437984497
$ws.Range("J7").Value = 0.04197179298073423
$ws.Range("K7").Value = 0.8470006184291898
$ws.Range("L7").Value = 0.04364069098483705

# Row 8 - SVM
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8686342689568496
$ws.Range("D8").Value = 0.02169271019065537
$ws.Range("E8").Value = 0.8787277072991359
$ws.Range("F8").Value = 0.03135104519053251
$ws.Range("G8").Value = 0.8791579262810908
$ws.Range("H8").Value = 0.02128888208214228
$ws.Range("I8").Value = 0.8811107073643412
$ws.Range("J8").Value = 0.03136849859284919
$ws.Range("K8").Value = 0.8612347969490827
$ws.Range("L8").Value = 0.0379998346360849

# Carry the existing header formatting (bold, centered, bordered) from
# B1 across the newly added header cells C1:L1.
$ws.Range("B1").Copy()
$ws.Range("C1:L1").PasteSpecial(-4122)

# Column A keeps its existing bold/bordered "index" styling; make sure
# the carried-down rows match (A2 already has it after the row delete).
$ws.Range("A2").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)

$excel.CutCopyMode = 0
